$d = $word.ActiveDocument
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$tbl = $ftr.Range.Tables(1)
$cell = $tbl.Cell(1,1)
$r = $cell.Range.Duplicate
$found = $r.Find.Execute("SocialManager", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found=$found Start=$($r.Start) End=$($r.End) Text=[$($r.Text)]"
